$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Activation costs added to optimization cost calculation:
# fuel_price (column F) for the "Hydro"/area-3 plant rows (17-40) is
# scaled up by a factor of 3.6 to account for activation costs.
$ws.Range("F17").Value = 2.0160000000000005
$ws.Range("F18").Value = 4.464
$ws.Range("F19").Value = 6.768000000000001
$ws.Range("F20").Value = 10.656
$ws.Range("F21").Value = 10.368
$ws.Range("F22").Value = 11.808000000000002
$ws.Range("F23").Value = 5.436
$ws.Range("F24").Value = 11.808000000000002
$ws.Range("F25").Value = 27.144000000000002
$ws.Range("F26").Value = 35.892
$ws.Range("F27").Value = 40.176
$ws.Range("F28").Value = 35.892
$ws.Range("F29").Value = 27.144000000000002
$ws.Range("F30").Value = 20.124000000000002
$ws.Range("F31").Value = 16.704
$ws.Range("F32").Value = 22.572000000000003
$ws.Range("F33").Value = 38.412000000000006
$ws.Range("F34").Value = 52.56
$ws.Range("F35").Value = 56.844
$ws.Range("F36").Value = 65.16
$ws.Range("F37").Value = 40.89600000000001
$ws.Range("F38").Value = 38.412000000000006
$ws.Range("F39").Value = 40.89600000000001
$ws.Range("F40").Value = 38.412000000000006
